$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D75:W106").Value = 0

$ws.Range("B84").Select()
